$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 2, shifting the existing row 2 down to row 3
$ws.Rows("2:2").Insert()
$ws.Rows("2:2").ClearFormats()

# Populate the newly inserted row 2 with the new match data
$ws.Range("A2").Value = 'hCSr7pE1'
$ws.Range("B2").Value = '18/11/2024'
$ws.Range("C2").Value = '17:00'
$ws.Range("D2").Value = 'ARGENTINA - TORNEO BETANO'
$ws.Range("E2").Value = 'Banfield'
$ws.Range("F2").Value = 'Tigre'
$ws.Range("G2").Value = 2.25
$ws.Range("H2").Value = 3.25
$ws.Range("I2").Value = 3.25
$ws.Range("J2").Value = 3.1
$ws.Range("K2").Value = 1.91
$ws.Range("L2").Value = 4.33
$ws.Range("M2").Value = 1.1
$ws.Range("N2").Value = 7
$ws.Range("O2").Value = 1.5
$ws.Range("P2").Value = 2.5
$ws.Range("Q2").Value = 2.6
$ws.Range("R2").Value = 1.48
$ws.Range("S2").Value = 1.57
$ws.Range("T2").Value = 2.25
$ws.Range("U2").Value = 2.1
$ws.Range("V2").Value = 1.67
$ws.Range("W2").Value = 6
$ws.Range("X2").Value = 9
$ws.Range("Y2").Value = 10
$ws.Range("Z2").Value = 21
$ws.Range("AA2").Value = 21
$ws.Range("AB2").Value = 41
$ws.Range("AC2").Value = 7
$ws.Range("AD2").Value = 6.5
$ws.Range("AE2").Value = 19
$ws.Range("AF2").Value = 81
$ws.Range("AG2").Value = 7.5
$ws.Range("AH2").Value = 15
$ws.Range("AI2").Value = 13
$ws.Range("AJ2").Value = 41
$ws.Range("AK2").Value = 34
$ws.Range("AL2").Value = 41
$ws.Range("AM2").Value = 201
$ws.Range("AN2").Value = 4
$ws.Range("AO2").Value = 13
$ws.Range("AP2").Value = 29
$ws.Range("AQ2").Value = 51
$ws.Range("AR2").Value = 81
$ws.Range("AS2").Value = 301
$ws.Range("AT2").Value = 2.25
$ws.Range("AU2").Value = 9
$ws.Range("AV2").Value = 81
$ws.Range("AW2").Value = 5
$ws.Range("AX2").Value = 21
$ws.Range("AY2").Value = 34
$ws.Range("AZ2").Value = 81
$ws.Range("BA2").Value = 126
$ws.Range("BB2").Value = 351
$ws.Range("BC2").Value = 126
$ws.Range("BD2").Value = 151

# Reorder the Correct-Score header labels AG1:AM1 (move "Odd_CS_4-4" from the front to the back)
$ws.Range("AG1").Value = 'Odd_CS_0-1'
$ws.Range("AH1").Value = 'Odd_CS_0-2'
$ws.Range("AI1").Value = 'Odd_CS_1-2'
$ws.Range("AJ1").Value = 'Odd_CS_0-3'
$ws.Range("AK1").Value = 'Odd_CS_1-3'
$ws.Range("AL1").Value = 'Odd_CS_2-3'
$ws.Range("AM1").Value = 'Odd_CS_4-4'

# Apply the same reorder to the shifted data row (row 3, formerly row 2)
$ws.Range("AG3").Value = 8
$ws.Range("AH3").Value = 13
$ws.Range("AI3").Value = 11
$ws.Range("AJ3").Value = 29
$ws.Range("AK3").Value = 23
$ws.Range("AL3").Value = 34
$ws.Range("AM3").Value = 351